# le15_p21_calls.xlsx — add an "Rat" (per-animal ID) column, populate it per
# source file, refresh the filter range to include it, and turn on AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for the new column -------------------------------------------
$ws.Range("T1").Value = "Rat"

# --- Per-animal "Rat" ID values, one block per source .mat file ----------
# (rat1_iso2, rat1_iso1) -> 563, (rat2_iso1, rat2_iso2) -> 566,
# (rat3_iso1, rat3_iso2) -> 565, (rat4_iso1, rat4_iso2) -> 564,
# (rat5_iso1, rat5_iso2) -> 568
$ws.Range("T2:T15").Value = 563
$ws.Range("T16:T19").Value = 563
$ws.Range("T20:T31").Value = 566
$ws.Range("T32:T56").Value = 566
$ws.Range("T57:T70").Value = 565
$ws.Range("T71:T85").Value = 565
$ws.Range("T86:T113").Value = 564
$ws.Range("T114:T134").Value = 564
$ws.Range("T135:T157").Value = 568
$ws.Range("T158:T168").Value = 568

# --- Extend the filter database / turn AutoFilter on over the new range --
$ws.Range("A1:T168").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$T`$168"
    }
}

# --- Match the author's final selection / scroll state -------------------
$ws.Range("T136:T168").Select() | Out-Null
